$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new columns before column H (old H:AE data shifts to K:AH)
$ws.Range("H1:J1").EntireColumn.Insert()

# New header row cells (H1, I1, J1)
$ws.Range("H1").Value = "headers"
$ws.Range("I1").Value = "Brushes"
$ws.Range("J1").Value = "Detangling Brushes"

# New row 18 - PLP category/breadcrumb data
$ws.Range("A18").Value = "Hair Tools"
$ws.Range("H18").Value = "Hair Tools"
$ws.Range("I18").Value = "Brushes"
$ws.Range("J18").Value = "Detangling Brushes"

# Blank cells with quote-prefix style (copy format from an existing cell using that style)
$ws.Range("R2").Copy()
$ws.Range("P18:R18").PasteSpecial(-4122)

# Blank cells needing the new black-font style
$ws.Range("V18:Y18").Font.Color = 0

# New row 19 - PLP product data
$ws.Range("A19").Value = "PLP Product"
$ws.Range("X19").Value = "The Lemon Bar Paddle Brush"
$ws.Range("Y19").Value = "1"
$ws.Range("R2").Copy()
$ws.Range("X19:Y19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Fix hyperlinks so they point at the shifted cells
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("K16"), "mailto:avayugundla@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("K17"), "mailto:avayugundla@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lotuswave@123")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lotuswave@123")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:avayugundla@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:avayugundla@helenoftroy.com")

# Update selection to match target workbook state
$ws.Range("A19").Select()
